$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure text date from 2021-05-21 to 2021-05-24
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-24 for illustrative purposes only and are subject to change."
$ws.Rows.Item(9).AutoFit()

# Update weight (D) and percent change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2552504219766386
$ws.Range("E2").Value = 0.008891454965358037

$ws.Range("D3").Value = 0.2556347780639556
$ws.Range("E3").Value = 0.004789781798829162

$ws.Range("D4").Value = 0.2434848758288402
$ws.Range("E4").Value = 0.006211180124223725

$ws.Range("D5").Value = 0.2456299241305656
$ws.Range("E5").Value = 0.01821493624772308

$ws.Range("D6").Value = 0.9999999999999999
$ws.Range("E6").Value = 0.00948044426888095

# Restore sheet protection (content/objects/scenarios locked; row & column
# formatting left unrestricted), matching the original workbook's settings.
$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
